$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was M, now B)
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9161290322580645
$ws.Range("C2").Value = 0.993006993006993
$ws.Range("D2").Value = 0.9530201342281879
$ws.Range("E2").Value = 143

# Row 3 (was B, now M)
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.9863013698630136
$ws.Range("C3").Value = 0.8470588235294118
$ws.Range("D3").Value = 0.9113924050632911
$ws.Range("E3").Value = 85

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.9385964912280702
$ws.Range("C4").Value = 0.9385964912280702
$ws.Range("D4").Value = 0.9385964912280702
$ws.Range("E4").Value = 0.9385964912280702

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.951215201060539
$ws.Range("C5").Value = 0.9200329082682024
$ws.Range("D5").Value = 0.9322062696457395
$ws.Range("E5").Value = 228

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9422897721546465
$ws.Range("C6").Value = 0.9385964912280702
$ws.Range("D6").Value = 0.9375010246710992
$ws.Range("E6").Value = 228
